$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M)
$ws.Range("D:E").Insert(-4161)

# For each data row, copy number-format/style from column F (shifted data) into the
# two newly inserted blank columns D and E, then populate their values.
$ws.Cells.Item(7, 6).Copy() | Out-Null
$ws.Range("D7:E7").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 6).Copy() | Out-Null
$ws.Range("D8:E8").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 4).Value = 2908800
$ws.Cells.Item(8, 5).Value = 2889100
$ws.Cells.Item(9, 6).Copy() | Out-Null
$ws.Range("D9:E9").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(9, 4).Value = 2476700
$ws.Cells.Item(9, 5).Value = 2454000
$ws.Cells.Item(10, 6).Copy() | Out-Null
$ws.Range("D10:E10").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10, 4).Value = 432100
$ws.Cells.Item(10, 5).Value = 435100
$ws.Cells.Item(11, 6).Copy() | Out-Null
$ws.Range("D11:E11").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(12, 6).Copy() | Out-Null
$ws.Range("D12:E12").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 6).Copy() | Out-Null
$ws.Range("D13:E13").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 6).Copy() | Out-Null
$ws.Range("D14:E14").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(14, 4).Value = 16500
$ws.Cells.Item(14, 5).Value = 23200
$ws.Cells.Item(15, 6).Copy() | Out-Null
$ws.Range("D15:E15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 4).Value = 17100
$ws.Cells.Item(15, 5).Value = 17000
$ws.Cells.Item(16, 6).Copy() | Out-Null
$ws.Range("D16:E16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 6).Copy() | Out-Null
$ws.Range("D17:E17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 4).Value = 2834200
$ws.Cells.Item(17, 5).Value = 2810900
$ws.Cells.Item(18, 6).Copy() | Out-Null
$ws.Range("D18:E18").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(18, 4).Value = 74600
$ws.Cells.Item(18, 5).Value = 78200
$ws.Cells.Item(19, 6).Copy() | Out-Null
$ws.Range("D19:E19").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20, 6).Copy() | Out-Null
$ws.Range("D20:E20").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20, 4).Value = -600
$ws.Cells.Item(20, 5).Value = -1200
$ws.Cells.Item(21, 6).Copy() | Out-Null
$ws.Range("D21:E21").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(21, 4).Value = 91100
$ws.Cells.Item(21, 5).Value = 94100
$ws.Cells.Item(22, 6).Copy() | Out-Null
$ws.Range("D22:E22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22, 4).Value = 34300
$ws.Cells.Item(22, 5).Value = 32700
$ws.Cells.Item(23, 6).Copy() | Out-Null
$ws.Range("D23:E23").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23, 4).Value = 39700
$ws.Cells.Item(23, 5).Value = 44400
$ws.Cells.Item(24, 6).Copy() | Out-Null
$ws.Range("D24:E24").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(24, 4).Value = 8800
$ws.Cells.Item(24, 5).Value = 10300
$ws.Cells.Item(25, 6).Copy() | Out-Null
$ws.Range("D25:E25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 6).Copy() | Out-Null
$ws.Range("D26:E26").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(26, 4).Value = 30800
$ws.Cells.Item(26, 5).Value = 34100
$ws.Cells.Item(27, 6).Copy() | Out-Null
$ws.Range("D27:E27").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27, 4).Value = 29700
$ws.Cells.Item(27, 5).Value = 32900
$ws.Cells.Item(28, 6).Copy() | Out-Null
$ws.Range("D28:E28").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 6).Copy() | Out-Null
$ws.Range("D29:E29").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(29, 4).Value = -100
$ws.Cells.Item(29, 5).Value = 700
$ws.Cells.Item(30, 6).Copy() | Out-Null
$ws.Range("D30:E30").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 6).Copy() | Out-Null
$ws.Range("D31:E31").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 6).Copy() | Out-Null
$ws.Range("D32:E32").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(32, 4).Value = 600
$ws.Cells.Item(32, 5).Value = 1200
$ws.Cells.Item(33, 6).Copy() | Out-Null
$ws.Range("D33:E33").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(33, 4).Value = 29600
$ws.Cells.Item(33, 5).Value = 33600
$ws.Cells.Item(34, 6).Copy() | Out-Null
$ws.Range("D34:E34").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 6).Copy() | Out-Null
$ws.Range("D35:E35").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(35, 4).Value = 29600
$ws.Cells.Item(35, 5).Value = 33600
$ws.Cells.Item(38, 6).Copy() | Out-Null
$ws.Range("D38:E38").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(39, 6).Copy() | Out-Null
$ws.Range("D39:E39").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(40, 6).Copy() | Out-Null
$ws.Range("D40:E40").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(41, 6).Copy() | Out-Null
$ws.Range("D41:E41").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(41, 4).Value = 15900
$ws.Cells.Item(41, 5).Value = 32000
$ws.Cells.Item(42, 6).Copy() | Out-Null
$ws.Range("D42:E42").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 6).Copy() | Out-Null
$ws.Range("D43:E43").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(43, 4).Value = 459600
$ws.Cells.Item(43, 5).Value = 404900
$ws.Cells.Item(44, 6).Copy() | Out-Null
$ws.Range("D44:E44").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(44, 4).Value = 1844100
$ws.Cells.Item(44, 5).Value = 1733800
$ws.Cells.Item(45, 6).Copy() | Out-Null
$ws.Range("D45:E45").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(45, 4).Value = 82700
$ws.Cells.Item(45, 5).Value = 78000
$ws.Cells.Item(46, 6).Copy() | Out-Null
$ws.Range("D46:E46").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(46, 4).Value = 2402400
$ws.Cells.Item(46, 5).Value = 2248700
$ws.Cells.Item(47, 6).Copy() | Out-Null
$ws.Range("D47:E47").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 6).Copy() | Out-Null
$ws.Range("D48:E48").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(48, 4).Value = 1347800
$ws.Cells.Item(48, 5).Value = 1350900
$ws.Cells.Item(49, 6).Copy() | Out-Null
$ws.Range("D49:E49").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(49, 4).Value = 1223600
$ws.Cells.Item(49, 5).Value = 1245100
$ws.Cells.Item(50, 6).Copy() | Out-Null
$ws.Range("D50:E50").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 6).Copy() | Out-Null
$ws.Range("D51:E51").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 6).Copy() | Out-Null
$ws.Range("D52:E52").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(52, 4).Value = 27300
$ws.Cells.Item(52, 5).Value = 36200
$ws.Cells.Item(53, 6).Copy() | Out-Null
$ws.Range("D53:E53").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 6).Copy() | Out-Null
$ws.Range("D54:E54").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(54, 4).Value = 5001100
$ws.Cells.Item(54, 5).Value = 4880900
$ws.Cells.Item(55, 6).Copy() | Out-Null
$ws.Range("D55:E55").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(56, 6).Copy() | Out-Null
$ws.Range("D56:E56").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(57, 6).Copy() | Out-Null
$ws.Range("D57:E57").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(57, 4).Value = 419400
$ws.Cells.Item(57, 5).Value = 428400
$ws.Cells.Item(58, 6).Copy() | Out-Null
$ws.Range("D58:E58").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(58, 4).Value = 1769600
$ws.Cells.Item(58, 5).Value = 1554800
$ws.Cells.Item(59, 6).Copy() | Out-Null
$ws.Range("D59:E59").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(59, 4).Value = 197600
$ws.Cells.Item(59, 5).Value = 207400
$ws.Cells.Item(60, 6).Copy() | Out-Null
$ws.Range("D60:E60").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(60, 4).Value = 2386600
$ws.Cells.Item(60, 5).Value = 2190700
$ws.Cells.Item(61, 6).Copy() | Out-Null
$ws.Range("D61:E61").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(61, 4).Value = 1281500
$ws.Cells.Item(61, 5).Value = 1304100
$ws.Cells.Item(62, 6).Copy() | Out-Null
$ws.Range("D62:E62").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(62, 4).Value = 237300
$ws.Cells.Item(62, 5).Value = 237700
$ws.Cells.Item(63, 6).Copy() | Out-Null
$ws.Range("D63:E63").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 6).Copy() | Out-Null
$ws.Range("D64:E64").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 6).Copy() | Out-Null
$ws.Range("D65:E65").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 6).Copy() | Out-Null
$ws.Range("D66:E66").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(66, 4).Value = 3905400
$ws.Cells.Item(66, 5).Value = 3732500
$ws.Cells.Item(67, 6).Copy() | Out-Null
$ws.Range("D67:E67").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(68, 6).Copy() | Out-Null
$ws.Range("D68:E68").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 6).Copy() | Out-Null
$ws.Range("D69:E69").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 6).Copy() | Out-Null
$ws.Range("D70:E70").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 6).Copy() | Out-Null
$ws.Range("D71:E71").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 6).Copy() | Out-Null
$ws.Range("D72:E72").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(72, 4).Value = 1394800
$ws.Cells.Item(72, 5).Value = 1368900
$ws.Cells.Item(73, 6).Copy() | Out-Null
$ws.Range("D73:E73").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 6).Copy() | Out-Null
$ws.Range("D74:E74").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 6).Copy() | Out-Null
$ws.Range("D75:E75").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 6).Copy() | Out-Null
$ws.Range("D76:E76").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(76, 4).Value = 1095700
$ws.Cells.Item(76, 5).Value = 1148300
$ws.Cells.Item(77, 6).Copy() | Out-Null
$ws.Range("D77:E77").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 6).Copy() | Out-Null
$ws.Range("D80:E80").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 6).Copy() | Out-Null
$ws.Range("D81:E81").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(81, 4).Value = 29600
$ws.Cells.Item(81, 5).Value = 33600
$ws.Cells.Item(82, 6).Copy() | Out-Null
$ws.Range("D82:E82").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(83, 6).Copy() | Out-Null
$ws.Range("D83:E83").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(83, 4).Value = 17100
$ws.Cells.Item(83, 5).Value = 17000
$ws.Cells.Item(84, 6).Copy() | Out-Null
$ws.Range("D84:E84").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 6).Copy() | Out-Null
$ws.Range("D85:E85").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 6).Copy() | Out-Null
$ws.Range("D86:E86").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 6).Copy() | Out-Null
$ws.Range("D87:E87").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 6).Copy() | Out-Null
$ws.Range("D88:E88").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 6).Copy() | Out-Null
$ws.Range("D89:E89").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(89, 4).Value = -87400
$ws.Cells.Item(89, 5).Value = 94300
$ws.Cells.Item(90, 6).Copy() | Out-Null
$ws.Range("D90:E90").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(91, 6).Copy() | Out-Null
$ws.Range("D91:E91").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(91, 4).Value = -22800
$ws.Cells.Item(91, 5).Value = -30000
$ws.Cells.Item(92, 6).Copy() | Out-Null
$ws.Range("D92:E92").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 6).Copy() | Out-Null
$ws.Range("D93:E93").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 6).Copy() | Out-Null
$ws.Range("D94:E94").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(94, 4).Value = -22500
$ws.Cells.Item(94, 5).Value = -57600
$ws.Cells.Item(95, 6).Copy() | Out-Null
$ws.Range("D95:E95").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(96, 6).Copy() | Out-Null
$ws.Range("D96:E96").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(96, 4).Value = -4900
$ws.Cells.Item(96, 5).Value = -5200
$ws.Cells.Item(97, 6).Copy() | Out-Null
$ws.Range("D97:E97").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 6).Copy() | Out-Null
$ws.Range("D98:E98").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 6).Copy() | Out-Null
$ws.Range("D99:E99").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 6).Copy() | Out-Null
$ws.Range("D100:E100").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(100, 4).Value = 94700
$ws.Cells.Item(100, 5).Value = -45800
$ws.Cells.Item(101, 6).Copy() | Out-Null
$ws.Range("D101:E101").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(101, 4).Value = -400
$ws.Cells.Item(101, 5).Value = -100
$ws.Cells.Item(102, 6).Copy() | Out-Null
$ws.Range("D102:E102").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(102, 4).Value = -15700
$ws.Cells.Item(102, 5).Value = -9200

$excel.CutCopyMode = 0

# A few historical quarters were restated by the author; update those specific cells
# in the columns that now hold the previously-existing data (post shift).
$ws.Cells.Item(61, 6).Value = 1358000
$ws.Cells.Item(61, 8).Value = 1318200
$ws.Cells.Item(62, 6).Value = 239400
$ws.Cells.Item(62, 8).Value = 230100
$ws.Cells.Item(91, 9).Value = -77000
$ws.Cells.Item(91, 10).Value = -67300
